$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.413.99"
$ws.Range("E2").Value = "  +0.78%  "

$ws.Range("D3").Value = "'1.839.93"
$ws.Range("E3").Value = "  +3.42%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'225.11"
$ws.Range("E5").Value = "  -0.14%  "

$ws.Range("D6").Value = "'0.559"
$ws.Range("E6").Value = "  +1.78%  "

$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "'32.11"
$ws.Range("E8").Value = "  +1.42%  "

$ws.Range("E9").Value = "  +4.27%  "

$ws.Range("D10").Value = "'0.0715"
$ws.Range("E10").Value = "  +8.92%  "

$ws.Range("E11").Value = "  +0.51%  "

$ws.Range("D12").Value = "'2.103.38"
$ws.Range("E12").Value = "  +3.39%  "

$ws.Range("D13").Value = "'1.845.33"
$ws.Range("E13").Value = "  +3.84%  "

$ws.Range("D14").Value = "'10.88"
$ws.Range("E14").Value = "  -1.88%  "

$ws.Range("D15").Value = "'0.649"
$ws.Range("E15").Value = "  +4.00%  "

$ws.Range("D16").Value = "'34.442.29"
$ws.Range("E16").Value = "  +0.94%  "

$ws.Range("D17").Value = "'4.37"
$ws.Range("E17").Value = "  +3.82%  "

$ws.Range("D18").Value = "'69.89"
$ws.Range("E18").Value = "  +1.70%  "

$ws.Range("D19").Value = "'252.11"
$ws.Range("E19").Value = "  -0.75%  "

$ws.Range("E20").Value = "  +8.24%  "

$ws.Range("D21").Value = "'11.41"
$ws.Range("E21").Value = "  +10.25%  "

$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("D23").Value = "'4.30"
$ws.Range("E23").Value = "  +2.47%  "

$ws.Range("E24").Value = "  +1.47%  "

$ws.Range("D25").Value = "'160.58"
$ws.Range("E25").Value = "  +2.49%  "

$ws.Range("D26").Value = "'16.74"
$ws.Range("E26").Value = "  +2.07%  "

$ws.Range("D27").Value = "'7.30"
$ws.Range("E27").Value = "  +4.44%  "

$ws.Range("E28").Value = "  +2.19%  "

$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("D30").Value = "'0.0538"
$ws.Range("E30").Value = "  +5.20%  "

$ws.Range("D31").Value = "'3.83"
$ws.Range("E31").Value = "  +1.55%  "

$ws.Range("E32").Value = "  +1.61%  "

$ws.Range("E33").Value = "  +1.42%  "

$ws.Range("E34").Value = "  +4.23%  "

$ws.Range("D35").Value = "'1.458.07"
$ws.Range("E35").Value = "  +1.24%  "

$ws.Range("D36").Value = "'0.649"
$ws.Range("E36").Value = "  +3.98%  "

$ws.Range("E37").Value = "  +1.70%  "

$ws.Range("E38").Value = "  +3.12%  "

$ws.Range("D39").Value = "'0.970"
$ws.Range("E39").Value = "  +9.15%  "

$ws.Range("D40").Value = "'82.47"
$ws.Range("E40").Value = "  -0.40%  "

$ws.Range("E41").Value = "  -2.53%  "

$ws.Range("E42").Value = "  +0.43%  "

$ws.Range("D43").Value = "'2.16"
$ws.Range("E43").Value = "  +5.32%  "

$ws.Range("E44").Value = "  +5.27%  "

$ws.Range("D45").Value = "'2.001.88"
$ws.Range("E45").Value = "  +3.34%  "

$ws.Range("D46").Value = "'0.0501"
$ws.Range("E46").Value = "  -1.68%  "

$ws.Range("E47").Value = "  +0.74%  "

$ws.Range("D48").Value = "'106.67"
$ws.Range("E48").Value = "  +8.42%  "

$ws.Range("D49").Value = "'12.13"
$ws.Range("E49").Value = "  +1.26%  "

$ws.Range("D50").Value = "'0.998"
$ws.Range("E50").Value = "  -0.04%  "

$ws.Range("D51").Value = "0.0₆0126"
$ws.Range("E51").Value = "  +8.40%  "

